$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing M-COD/PRODUCTO/...
# columns one place to the right (A->B, B->C, ... F->G)
$ws.Range("A1").EntireColumn.Insert()

# New header cell for the inserted column
$ws.Range("A1").Value = "FILTRO"

# Copy the header formatting (bold, centered, bordered) from the neighboring
# header cell onto the new header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new column's data rows with the "NOTEBOOK" filter value
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = "NOTEBOOK"
}
